# Issue #15 remove dirname from path
# - Mark existing Issue #15 (row 16 on the "Issues" sheet) as DONE.
# - Log a new Issue #17 (row 18): "add return option to nav".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Issue #15 ("remove dirname from path") is now complete -> set Status to DONE.
$ws.Range("B16").Value = "DONE"

# Add new issue row: #17 "add return option to nav"
$ws.Range("A18").Value = 17
$ws.Range("D18").Value = "add return option to nav"

# Move the active selection to the newly added cell, matching the saved view state.
$ws.Range("D18").Select()
